# Regenerate the "K" column (column G) values in the save_data sheet.
# The underlying data generation changed from using "Strike#" to "K",
# so the raw strikeout counts recorded in column G need to be rewritten
# with the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$newValues = @{
    2  = 5
    3  = 6
    4  = 2
    5  = 6
    6  = 4
    7  = 3
    8  = 6
    9  = 9
    10 = 1
    11 = 6
    12 = 11
    13 = 3
    14 = 4
    15 = 5
    16 = 4
    17 = 3
    18 = 0
    19 = 2
    20 = 2
    21 = 5
    22 = 4
    23 = 6
    24 = 5
    25 = 3
    26 = 4
    27 = 6
    28 = 2
    29 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
